# The deck's live theme (ppt/theme/theme2.xml - the one actually wired to
# the slide master via presentation.xml.rels / slideMaster1.xml.rels, and
# therefore the one every slide/layout renders with) currently carries the
# "Integral" palette. The target revision swaps it back to the plain
# default "Office Theme" palette. Drive this through the Slide's
# ThemeColorScheme (PowerPoint's 12-slot DrawingML scheme, in order: dk1,
# lt1, dk2, lt2, accent1-6, hlink, folHlink) - writing through any slide
# updates the single shared master theme that every slide/layout inherits.
#
# RGB() isn't available in this host, so each value below is the decimal
# form of the target hex color (0xBBGGRR order, i.e. R + G*256 + B*65536).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      #000000
$tcs.Item(2).RGB  = 16777215   # lt1      #FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      #44546A
$tcs.Item(4).RGB  = 15132391   # lt2      #E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  #5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  #ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  #A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  #FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  #4472C4
$tcs.Item(10).RGB = 4697456    # accent6  #70AD47
$tcs.Item(11).RGB = 12673797   # hlink    #0563C1
$tcs.Item(12).RGB = 7491477    # folHlink #954F72
